$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: ABN number changed, Entity name changed (merged shared strings),
# Second Business Name now reads "N/A" (same as Business Name).
$ws.Range("B3").Value = 95488716489
$ws.Range("C3").Value = "I.C HILL & L.A HILL & S HILL & R.L JAKINS"
$ws.Range("D3").Value = "N/A"
$ws.Range("E3").Value = "N/A"

# Row 5: Business Name / Second Business Name now read "N/A" as well
# (unaffected values kept, just re-asserted through shared-string reindex).
$ws.Range("D5").Value = "N/A"
$ws.Range("E5").Value = "N/A"

# Update the active selection to C4 (was C5).
$ws.Range("C4").Select()
